$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.983.77"
$ws.Range("E2").Value = "  +4.06%  "

$ws.Range("D3").Value = "2.480.17"
$ws.Range("E3").Value = "  +2.09%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'491.55"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").Value = "'151.62"
$ws.Range("E6").Value = "  +9.26%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +2.83%  "

$ws.Range("D9").Value = "2.491.83"
$ws.Range("E9").Value = "  +1.64%  "

$ws.Range("D10").Value = "'0.0999"
$ws.Range("E10").Value = "  +4.31%  "

$ws.Range("D11").Value = "'5.73"
$ws.Range("E11").Value = "  +4.87%  "

$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  +4.22%  "

$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").Value = "2.916.78"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "57.200.12"
$ws.Range("E15").Value = "  +4.23%  "

$ws.Range("D16").Value = "'21.08"
$ws.Range("E16").Value = "  +3.02%  "

$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "2.491.88"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").Value = "'4.55"
$ws.Range("E19").Value = "  +5.17%  "

$ws.Range("D20").Value = "'10.22"
$ws.Range("E20").Value = "  +3.99%  "

$ws.Range("D21").Value = "'320.90"
$ws.Range("E21").Value = "  +2.74%  "

$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'5.90"
$ws.Range("E23").Value = "  +4.53%  "

$ws.Range("D24").Value = "'58.22"
$ws.Range("E24").Value = "  +1.94%  "

$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "'0.407"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.164"
$ws.Range("E27").Value = "  +2.29%  "

$ws.Range("D28").Value = "2.617.49"
$ws.Range("E28").Value = "  +2.41%  "

$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  +2.71%  "

$ws.Range("D30").Value = "0.0₃0812"
$ws.Range("E30").Value = "  +5.83%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").Value = "'151.41"
$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("D33").Value = "'18.31"
$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("D34").Value = "'1.52"
$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("D35").Value = "'5.23"
$ws.Range("E35").Value = "  +1.68%  "

$ws.Range("D36").Value = "'0.896"
$ws.Range("E36").Value = "  +5.84%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.15"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.78"
$ws.Range("E38").Value = "  +5.50%  "

$ws.Range("D39").Value = "'1.40"
$ws.Range("E39").Value = "  +8.79%  "

$ws.Range("D40").Value = "'34.19"
$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = "  +3.09%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0558"
$ws.Range("E42").Value = "  +3.10%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.612"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").Value = "'4.85"
$ws.Range("E45").Value = "  +4.52%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0942"
$ws.Range("E46").Value = "  +6.43%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'261.88"
$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.22"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = "  +3.25%  "

$ws.Range("D50").Value = "'17.82"
$ws.Range("E50").Value = "  +4.41%  "

$ws.Range("D51").Value = "1.875.78"
$ws.Range("E51").Value = "  -2.24%  "

